$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '61.999.54'
$ws.Range("E2").Value = '  +0.76%  '
# Row 3
$ws.Range("D3").Value = '3.416.97'
$ws.Range("E3").Value = '  +1.12%  '
# Row 4
$ws.Range("E4").Value = '  -0.09%  '
# Row 5
$ws.Range("D5").Value = '''409.63'
$ws.Range("E5").Value = '  +0.71%  '
# Row 6
$ws.Range("D6").Value = '''128.87'
$ws.Range("E6").Value = '  -4.09%  '
# Row 7
$ws.Range("D7").Value = '''0.625'
$ws.Range("E7").Value = '  +5.33%  '
# Row 8
$ws.Range("E8").Value = '  -0.12%  '
# Row 9
$ws.Range("D9").Value = '''0.749'
$ws.Range("E9").Value = '  +11.70%  '
# Row 10
$ws.Range("D10").Value = '''0.141'
$ws.Range("E10").Value = '  +16.79%  '
# Row 11
$ws.Range("D11").Value = '''42.96'
$ws.Range("E11").Value = '  +0.78%  '
# Row 13
$ws.Range("D13").Value = '''21.17'
$ws.Range("E13").Value = '  +7.39%  '
# Row 14
$ws.Range("D14").Value = '''8.87'
$ws.Range("E14").Value = '  +5.30%  '
# Row 15
$ws.Range("D15").Value = '''0.0000202'
$ws.Range("E15").Value = '  +58.79%  '
# Row 16
$ws.Range("D16").Value = '3.417.28'
$ws.Range("E16").Value = '  +0.17%  '
# Row 17
$ws.Range("D17").Value = '''12.63'
$ws.Range("E17").Value = '  +14.72%  '
# Row 18
$ws.Range("E18").Value = '  +3.82%  '
# Row 19
$ws.Range("D19").Value = '61.984.18'
$ws.Range("E19").Value = '  +0.87%  '
# Row 20
$ws.Range("D20").Value = '''404.47'
$ws.Range("E20").Value = '  +28.47%  '
# Row 21
$ws.Range("D21").Value = '''90.46'
$ws.Range("E21").Value = '  +6.29%  '
# Row 22
$ws.Range("E22").Value = '  -0.79%  '
# Row 23
$ws.Range("D23").Value = '''13.44'
$ws.Range("E23").Value = '  +4.80%  '
# Row 24
$ws.Range("D24").Value = '''3.24'
$ws.Range("E24").Value = '  +3.24%  '
# Row 25
$ws.Range("D25").Value = '''32.92'
$ws.Range("E25").Value = '  +11.28%  '
# Row 26
$ws.Range("D26").Value = '''4.79'
$ws.Range("E26").Value = '  -0.03%  '
# Row 27
$ws.Range("D27").Value = '''8.52'
$ws.Range("E27").Value = '  +1.77%  '
# Row 28
$ws.Range("D28").Value = '''7.64'
$ws.Range("E28").Value = '  +0.21%  '
# Row 29
$ws.Range("D29").Value = '''2.72'
$ws.Range("E29").Value = '  +4.21%  '
# Row 30
$ws.Range("E30").Value = '  +1.14%  '
# Row 31
$ws.Range("E31").Value = '  +0.62%  '
# Row 32
$ws.Range("D32").Value = '''43.99'
$ws.Range("E32").Value = '  +8.00%  '
# Row 33
$ws.Range("D33").Value = '''11.82'
$ws.Range("E33").Value = '  +4.05%  '
# Row 34
$ws.Range("D34").Value = '''1.00'
$ws.Range("E34").Value = '  +0.05%  '
# Row 35
$ws.Range("D35").Value = '''0.0501'
$ws.Range("E35").Value = '  +3.72%  '
# Row 36
$ws.Range("D36").Value = '''52.92'
$ws.Range("E36").Value = '  +2.01%  '
# Row 37
$ws.Range("E37").Value = '  +0.05%  '
# Row 38
$ws.Range("E38").Value = '  -0.79%  '
# Row 39
$ws.Range("E39").Value = '  +6.47%  '
# Row 40
$ws.Range("E40").Value = '  -1.08%  '
# Row 41
$ws.Range("E41").Value = '  +6.21%  '
# Row 42
$ws.Range("D42").Value = '''141.04'
$ws.Range("E42").Value = '  +1.41%  '
# Row 43
$ws.Range("E43").Value = '  -0.36%  '
# Row 44
$ws.Range("D44").Value = '''4.02'
$ws.Range("E44").Value = '  -0.53%  '
# Row 45
$ws.Range("D45").Value = '''16.80'
$ws.Range("E45").Value = '  +0.36%  '
# Row 46
$ws.Range("E46").Value = '  +6.01%  '
# Row 47
$ws.Range("D47").Value = '''21.85'
$ws.Range("E47").Value = '  +2.88%  '
# Row 48
$ws.Range("D48").Value = '2.108.90'
$ws.Range("E48").Value = '  -0.68%  '
# Row 49
$ws.Range("D49").Value = '''2.28'
$ws.Range("E49").Value = '  -0.61%  '
# Row 50
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '''0.127'
$ws.Range("E50").Value = '  +14.34%  '
# Row 51
$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D51").Value = '''1.93'
$ws.Range("E51").Value = '  +0.45%  '
